# Apply weekly report value updates to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (B, C, D) new values. Row 19 has no B value (string "---" stays as-is).
$updates = @{
    2  = @(0.52,  -22.46,              805.3855738830566)
    3  = @(-0.99, 17.6,                1325.573728547231)
    4  = @(1.74,  -34.19,              271.8245054743797)
    5  = @(-1.23, -31.14,              178.3013508430562)
    6  = @(-2.34, 152.75,              1436.618495895817)
    7  = @(4.85,  -7.13,               680.9286308111443)
    8  = @(5.4,   -21.03,              367.6451586415828)
    9  = @(3.67,  -64.54000000000001,  428.0999908447266)
    10 = @(-2.33, -0.4,                602.4)
    11 = @(-0.15, 4.85,                532.9424989723205)
    12 = @(-1.18, 36.22,               809.3524568252564)
    13 = @(-2.81, -20.01,              1201.99974420556)
    14 = @(4.16,  43.54,               1016.161875359902)
    15 = @(0.33,  -30.9,               485.3440982540862)
    16 = @(0.87,  34.03,               597.668439994812)
    17 = @(0.63,  -65.15000000000001,  134.4812603407771)
    18 = @(0.12,  8.75,                10874.72780889371)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
}

# Row 19: only C and D change; B19 remains the "---" string.
$ws.Cells.Item(19, 3).Value = 1346.635221111907
$ws.Cells.Item(19, 4).Value = 935.8946862732214
